$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (shifts existing rows 32-99 down to 33-100)
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with the new record
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44838
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112021
$ws.Range("G32").Value = "Ají"
$ws.Range("H32").Value = "Inferno"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 60
$ws.Range("K32").Value = 24000
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 24500
$ws.Range("N32").Value = "`$/caja 10 kilos"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 2450
$ws.Range("Q32").Value = 10
$ws.Range("R32").Value = "Hortaliza"

# Append a brand-new row 101 at the end of the table
$ws.Range("A101").Value = 7
$ws.Range("B101").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C101").Value = "Ñuble"
$ws.Range("D101").Value = 44832
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E101").Value = 16
$ws.Range("F101").Value = 100112021
$ws.Range("G101").Value = "Ají"
$ws.Range("H101").Value = "Inferno"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 60
$ws.Range("K101").Value = 24000
$ws.Range("L101").Value = 25000
$ws.Range("M101").Value = 24500
$ws.Range("N101").Value = "`$/caja 10 kilos"
$ws.Range("O101").Value = "Región de Arica y Parinacota"
$ws.Range("P101").Value = 2450
$ws.Range("Q101").Value = 10
$ws.Range("R101").Value = "Hortaliza"
